{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Paragraph 1: \"Tyler Roop\" -> \"Tyler \" run + spell-check-flagged \"Roop\" run\n// (mirrors Word's own proofing pass marking \"Roop\" as a misspelling after\n// the author edited the paragraph).\nconst firstPara = paragraphs.items[0];\nconst firstOoxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n          '<w:body>' +\n            '<w:p>' +\n              '<w:r><w:t xml:space=\"preserve\">Tyler </w:t></w:r>' +\n              '<w:proofErr w:type=\"spellStart\"/>' +\n              '<w:r><w:t>Roop</w:t></w:r>' +\n              '<w:proofErr w:type=\"spellEnd\"/>' +\n            '</w:p>' +\n          '</w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>';\nfirstPara.insertOoxml(firstOoxml, Word.InsertLocation.replace);\n\n// Paragraph 2: \"Weihan Huang\" -> \"Lance Grengbondai\" (plain text swap; the\n// bookmarkStart/bookmarkEnd around \"_GoBack\" stay put since we replace the\n// paragraph's own text run, not the whole paragraph).\nconst secondPara = paragraphs.items[1];\nsecondPara.insertText(\"Lance Grengbondai\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Paragraph 1: \"Tyler Roop\" -> \"Tyler \" run + spell-check-flagged \"Roop\" run\n# (mirrors Word's own proofing pass marking \"Roop\" as a misspelling after\n# the author edited the paragraph).\n$para1 = $d.Paragraphs(1).Range\n$ooxml = @\"\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t xml:space=\"preserve\">Tyler </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>Roop</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n$para1.InsertXML($ooxml)\n\n# Paragraph 2: \"Weihan Huang\" -> \"Lance Grengbondai\" (targeted find/replace so\n# the bookmarkStart/bookmarkEnd around \"_GoBack\" is left untouched).\n$find = $d.Content.Find\n$find.Text = \"Weihan Huang\"\n$find.Replacement.Text = \"Lance Grengbondai\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
